# Refresh the cryptos price/volume table (GitHub Actions style data sync).
# Note: Price column (D) values are stored as text (e.g. "43.308.26" uses
# '.' as a thousands separator, not a decimal point), so a leading
# apostrophe ('' inside a single-quoted PS string -> a literal ') is used
# to force Excel to keep them as text instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''43.308.26'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '''2.603.30'
$ws.Range('E3').Value = '  +3.81%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''316.71'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = '''98.40'
$ws.Range('E6').Value = '  +4.82%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.542'
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('D10').Value = '''36.14'
$ws.Range('E10').Value = '  +1.65%  '
$ws.Range('D11').Value = '''0.0816'
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('D12').Value = '''7.58'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').Value = '''3.004.12'
$ws.Range('E13').Value = '  +3.89%  '
$ws.Range('D15').Value = '''2.595.75'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').Value = '''15.26'
$ws.Range('E16').Value = '  +1.26%  '
$ws.Range('D17').Value = '''0.851'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '''43.437.35'
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('D19').Value = '''6.88'
$ws.Range('E19').Value = '  +3.97%  '
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').Value = '''0.0₃0972'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').Value = '''69.62'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '''254.93'
$ws.Range('E23').Value = '  +2.25%  '
$ws.Range('D24').Value = '''2.98'
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('E25').Value = '  +4.62%  '
$ws.Range('D26').Value = '''27.30'
$ws.Range('E26').Value = '  +2.69%  '
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').Value = '''2.43'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').Value = '''41.61'
$ws.Range('E29').Value = '  +3.24%  '
$ws.Range('D30').Value = '''10.34'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').Value = '''5.89'
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('D32').Value = '''156.84'
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('D33').Value = '''3.50'
$ws.Range('E33').Value = '  +7.14%  '
$ws.Range('D34').Value = '''2.17'
$ws.Range('E34').Value = '  +3.70%  '
$ws.Range('D35').Value = '''0.0813'
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('E36').Value = '  +3.02%  '
$ws.Range('D37').Value = '''18.86'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '''0.113'
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('B39').Value = 'ApeXProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D39').Value = '''2.51'
$ws.Range('E39').Value = '  +10.07%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = '''22.87'
$ws.Range('E41').Value = '  -2.11%  '
$ws.Range('D42').Value = '''4.02'
$ws.Range('E42').Value = '  +7.51%  '
$ws.Range('D43').Value = '''0.0306'
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D45').Value = '''3.26'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').Value = '''2.017.02'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').Value = '''2.858.23'
$ws.Range('E48').Value = '  +3.64%  '
$ws.Range('D49').Value = '''83.81'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').Value = '''75.18'
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('D51').Value = '''0.196'
$ws.Range('E51').Value = '  +4.97%  '
